$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write temp headers to a blank area (row 100) with default style
$headers = @(
  "Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210",
  "Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210",
  "diff",
  "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304",
  "Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304"
)
for ($i = 0; $i -lt $headers.Count; $i++) {
  $ws.Cells.Item(100, $i + 1).Value = $headers[$i]
}
$tmpRng = $ws.Range("A100:U100")
$tbl = $ws.ListObjects.Add(1, $tmpRng, $null, 1)
Write-Host "Created on temp range"
try {
  $tbl.Resize($ws.Range("A1:U57"))
  Write-Host "Resized ok"
} catch {
  Write-Host "resize err: $_"
}
